$d = $word.ActiveDocument

$replacements = @(
    @("30÷7=4, 2", "21÷8=2, 5"),
    @("50÷5=10, 0", "53÷6=8, 5"),
    @("29÷8=3, 5", "81÷2=40, 1"),
    @("95÷9=10, 5", "83÷4=20, 3"),
    @("77÷5=15, 2", "14÷4=3, 2"),
    @("99÷5=19, 4", "95÷3=31, 2"),
    @("10÷8=1, 2", "85÷9=9, 4"),
    @("10÷2=5, 0", "20÷9=2, 2"),
    @("45÷5=9, 0", "12÷7=1, 5"),
    @("17÷4=4, 1", "98÷8=12, 2"),
    @("24÷7=3, 3", "31÷3=10, 1"),
    @("24÷8=3, 0", "19÷7=2, 5"),
    @("31÷5=6, 1", "63÷4=15, 3"),
    @("24÷2=12, 0", "25÷2=12, 1"),
    @("14÷5=2, 4", "91÷6=15, 1"),
    @("70÷3=23, 1", "65÷8=8, 1"),
    @("52÷9=5, 7", "15÷9=1, 6"),
    @("76÷3=25, 1", "25÷2=12, 1"),
    @("91÷8=11, 3", "13÷9=1, 4"),
    @("83÷3=27, 2", "21÷3=7, 0"),
    @("77÷3=25, 2", "36÷8=4, 4"),
    @("64÷7=9, 1", "97÷2=48, 1"),
    @("74÷6=12, 2", "88÷2=44, 0"),
    @("57÷9=6, 3", "66÷8=8, 2"),
    @("17÷5=3, 2", "12÷4=3, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
